$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New block of rows describing an additional experiment run named
# "hh_n1_pre_az_test" (a rename/addition of the pre_az data file used by UJ).
# Columns: A = experiment, B = feature-selection method, C = feature count, D = model
$experiment = "hh_n1_pre_az_test"

$rows = @(
    @("pca", 5, "MLP"),
    @("pca", 5, "RandFor"),
    @("pca", 5, "SVC"),
    @("pca", 10, "MLP"),
    @("pca", 10, "RandFor"),
    @("pca", 10, "SVC"),
    @("pca", 15, "MLP"),
    @("pca", 15, "RandFor"),
    @("pca", 15, "SVC"),
    @("pca", 20, "MLP"),
    @("pca", 20, "RandFor"),
    @("pca", 20, "SVC"),
    @("pca", 25, "MLP"),
    @("pca", 25, "RandFor"),
    @("pca", 25, "SVC"),
    @("rf", 5, "MLP"),
    @("rf", 5, "RandFor"),
    @("rf", 5, "SVC"),
    @("rf", 10, "MLP"),
    @("rf", 10, "RandFor"),
    @("rf", 10, "SVC"),
    @("rf", 15, "MLP"),
    @("rf", 15, "RandFor"),
    @("rf", 15, "SVC"),
    @("rf", 20, "MLP"),
    @("rf", 20, "RandFor"),
    @("rf", 20, "SVC"),
    @("rf", 25, "MLP"),
    @("rf", 25, "RandFor"),
    @("rf", 25, "SVC"),
    @("svc", 5, "MLP"),
    @("svc", 5, "RandFor"),
    @("svc", 5, "SVC"),
    @("svc", 10, "MLP"),
    @("svc", 10, "RandFor"),
    @("svc", 10, "SVC"),
    @("svc", 15, "MLP"),
    @("svc", 15, "RandFor"),
    @("svc", 15, "SVC"),
    @("svc", 20, "MLP"),
    @("svc", 20, "RandFor"),
    @("svc", 20, "SVC"),
    @("svc", 25, "MLP"),
    @("svc", 25, "RandFor"),
    @("svc", 25, "SVC")
)

$startRow = 143
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $feat = $rows[$i][0]
    $cnt = $rows[$i][1]
    $model = $rows[$i][2]

    $ws.Cells.Item($r, 1).Value = $experiment
    $ws.Cells.Item($r, 2).Value = $feat
    $ws.Cells.Item($r, 3).Value = $cnt
    $ws.Cells.Item($r, 4).Value = $model
}

# Column A needs to widen to fit the new, longer experiment name.
$ws.Columns.Item(1).ColumnWidth = 17.7

# Restore the selection/view state similar to the authored file
$ws.Range("E143").Select()
